# The document's opening paragraph described the project in a single run.
# The edit reorders "can easily import" -> "easily can import" (the rest
# of the sentence/paragraph is unchanged). Use Find/Replace scoped to the
# unique surrounding phrase so only the intended occurrence is touched.

$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "owners can easily import",  # find what
    $true,                       # MatchCase
    $false,                      # MatchWholeWord
    $false,                      # MatchWildcards
    $false,                      # MatchSoundsLike
    $false,                      # MatchAllWordForms
    $true,                       # Forward
    1,                           # Wrap (wdFindContinue)
    $false,                      # Format
    "owners easily can import",  # replace with
    2                            # Replace (wdReplaceAll)
)

Write-Host "Replaced: $found"
Write-Host $d.Content.Text
